$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FORM IMUNISASI")

$ws.Range("A2").Value = "53.06.13.2001"
$ws.Range("A3").Value = "53.06.13.2002"
$ws.Range("A4").Value = "53.06.13.2003"
$ws.Range("A5").Value = "53.06.13.2004"
$ws.Range("A6").Value = "53.06.13.2005"
$ws.Range("A7").Value = "53.06.13.2006"
$ws.Range("A8").Value = "53.06.13.2007"
$ws.Range("A9").Value = "53.06.13.2008"

$ws.Range("A15").Select()
